# Update "想去人数" (F column) figures across the four sheets to match the
# newly generated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 32
$ws.Range("F6").Value  = 72
$ws.Range("F7").Value  = 287
$ws.Range("F8").Value  = 346
$ws.Range("F9").Value  = 3295
$ws.Range("F10").Value = 1168
$ws.Range("F11").Value = 1050
$ws.Range("F12").Value = 868
$ws.Range("F13").Value = 94
$ws.Range("F14").Value = 847
$ws.Range("F15").Value = 1499
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = 806
$ws.Range("F18").Value = 1710
$ws.Range("F20").Value = 381
$ws.Range("F21").Value = 188
$ws.Range("F22").Value = 74
$ws.Range("F23").Value = 117
$ws.Range("F25").Value = 2641

# --- Sheet "演出" -----------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F18").Value = 9
$ws.Range("F40").Value = 360
$ws.Range("F48").Value = 299

# --- Sheet "本地生活" -------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 2495
$ws.Range("F6").Value  = 2507
$ws.Range("F7").Value  = 9547
$ws.Range("F8").Value  = 141
$ws.Range("F11").Value = 364
$ws.Range("F12").Value = 2819
$ws.Range("F13").Value = 373
$ws.Range("F14").Value = 687

# --- Sheet "全部类型" -------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 2495
$ws.Range("F3").Value  = 141
$ws.Range("F7").Value  = 364
$ws.Range("F8").Value  = 373
$ws.Range("F13").Value = 72
$ws.Range("F14").Value = 287
$ws.Range("F15").Value = 346
$ws.Range("F16").Value = 1168
$ws.Range("F18").Value = 1050
$ws.Range("F19").Value = 868
$ws.Range("F20").Value = 94
$ws.Range("F21").Value = 847
$ws.Range("F23").Value = 1499
$ws.Range("F27").Value = 806
$ws.Range("F31").Value = 1710
$ws.Range("F32").Value = 381
$ws.Range("F39").Value = 74
$ws.Range("F44").Value = 360
$ws.Range("F45").Value = 2641
